# The deck currently has its Design/Theme color scheme set to the
# "Integral" palette (dk2/lt2/accent1-6/hlink/folHlink). This restores
# the stock Office Theme color palette, matching the Design > Colors
# swap recorded in the commit.

function Get-OleColor($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$cs = $s.ThemeColorScheme

# Theme color slots, in COM ThemeColorScheme.Item(n) order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$officeColors = @(
    (Get-OleColor 0x00 0x00 0x00), # dk1
    (Get-OleColor 0xFF 0xFF 0xFF), # lt1
    (Get-OleColor 0x44 0x54 0x6A), # dk2
    (Get-OleColor 0xE7 0xE6 0xE6), # lt2
    (Get-OleColor 0x5B 0x9B 0xD5), # accent1
    (Get-OleColor 0xED 0x7D 0x31), # accent2
    (Get-OleColor 0xA5 0xA5 0xA5), # accent3
    (Get-OleColor 0xFF 0xC0 0x00), # accent4
    (Get-OleColor 0x44 0x72 0xC4), # accent5
    (Get-OleColor 0x70 0xAD 0x47), # accent6
    (Get-OleColor 0x05 0x63 0xC1), # hlink
    (Get-OleColor 0x95 0x4F 0x72)  # folHlink
)

for ($i = 1; $i -le $cs.Count; $i++) {
    $cs.Item($i).RGB = $officeColors[$i - 1]
}
